# Apply "update script database xlsx" change to the r0 sheet:
# Add four new named script rows (ExoT_r0_script_8v0 .. 8v3) which correspond
# to the existing "Hanfei's flow rate optimization" rows (previously at rows
# 21-24, unnamed in column A) re-homed to rows 19-22 with their Name (column
# A) filled in.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("r0")

# Capture the values currently sitting in rows 21-24 (columns B:K) before we
# start overwriting anything.
$srcRows = @(21, 22, 23, 24)
$colLetters = @("B", "C", "D", "E", "F", "G", "H", "I", "J", "K")

$data = @{}
foreach ($r in $srcRows) {
    $rowVals = @{}
    foreach ($col in $colLetters) {
        $rowVals[$col] = $ws.Range("$col$r").Value2
    }
    $data[$r] = $rowVals
}

# Row height that the last of these rows (now row 22) ends up with.
$lastRowHeight = 16.95

# Destination rows 19-22, paired with the source rows that feed them, and the
# new script name for column A.
$mapping = @(
    @{ Dest = 19; Src = 21; Name = "ExoT_r0_script_8v0" },
    @{ Dest = 20; Src = 22; Name = "ExoT_r0_script_8v1" },
    @{ Dest = 21; Src = 23; Name = "ExoT_r0_script_8v2" },
    @{ Dest = 22; Src = 24; Name = "ExoT_r0_script_8v3" }
)

foreach ($m in $mapping) {
    $dest = $m.Dest
    $src = $m.Src
    $ws.Range("A$dest").Value = $m.Name
    foreach ($col in $colLetters) {
        $ws.Range("$col$dest").Value = $data[$src][$col]
    }
}

# The old row 22's custom height carries over to the new row 22.
$ws.Rows.Item(22).RowHeight = $lastRowHeight

# Rows 23 and 24 no longer exist in the final layout - clear them out
# entirely so the sheet dimension shrinks back down to row 22.
$ws.Rows("23:24").Delete()

# Freeze panes: split after column A and after row 3 (xSplit=1 / ySplit=3),
# matching the new frozen-pane view with the active cell in the bottom-right
# pane at G4.
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("B4").Select()
$excel.ActiveWindow.FreezePanes = $true

$ws.Range("G4").Select()
